$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.210.78'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '3.026.36'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.41'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.52'
$ws.Range("E6").Value = '  -2.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -2.05%  '
$ws.Range("D9").Value = '3.034.83'
$ws.Range("E9").Value = '  -1.82%  '
$ws.Range("E10").Value = '  -0.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  -3.96%  '
$ws.Range("E12").Value = '  -1.71%  '
$ws.Range("D13").Value = '3.555.85'
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("D15").Value = '63.255.49'
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.19'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000151'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '3.033.17'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '398.72'
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("E20").Value = '  +0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.04'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("E22").Value = '  -4.41%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").Value = '0.0₃0989'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.75'
$ws.Range("E28").Value = '  +2.61%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.51'
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.91'
$ws.Range("E33").Value = '  +4.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.76'
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("E35").Value = '  +2.72%  '
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.32'
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("D38").Value = '2.544.08'
$ws.Range("E38").Value = '  -5.73%  '
$ws.Range("E39").Value = '  -2.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.99'
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.96'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.61'
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0604'
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0250'
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.09'
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.998'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.24'
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '269.87'
$ws.Range("E49").Value = '  -2.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0952'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.50'
$ws.Range("E51").Value = '  +0.51%  '
